$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style from an untouched, default-styled cell (column B, coin names)
# used to restore cell style after forcing text entry for numeric-looking values,
# so we don't leave a stray NumberFormat behind on edited cells.
$cleanStyle = $ws.Range('B2').Style

$ws.Range('D2').Value = '26.403.87'
$ws.Range('E2').Value = '  +0.50%  '
$ws.Range('D3').Value = '1.695.80'
$ws.Range('E3').Value = '  +0.97%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.010'
$ws.Range('D4').Style = $cleanStyle
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.52'
$ws.Range('D5').Style = $cleanStyle
$ws.Range('E5').Value = '  +0.09%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5494'
$ws.Range('D6').Style = $cleanStyle
$ws.Range('E6').Value = '  +4.49%  '
$ws.Range('E7').Value = '  +0.19%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2735'
$ws.Range('D8').Style = $cleanStyle
$ws.Range('E8').Value = '  +1.39%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06444'
$ws.Range('D9').Style = $cleanStyle
$ws.Range('E9').Value = '  +0.16%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.95'
$ws.Range('D10').Style = $cleanStyle
$ws.Range('E10').Value = '  -0.21%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07668'
$ws.Range('D11').Style = $cleanStyle
$ws.Range('E11').Value = '  +2.21%  '
$ws.Range('D12').Value = '1.726.13'
$ws.Range('E12').Value = '  +2.42%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.556'
$ws.Range('D13').Style = $cleanStyle
$ws.Range('E13').Value = '  +0.30%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5845'
$ws.Range('D14').Style = $cleanStyle
$ws.Range('E14').Value = '  +0.67%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000008407'
$ws.Range('D15').Style = $cleanStyle
$ws.Range('E15').Value = '  -0.86%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.67'
$ws.Range('D16').Style = $cleanStyle
$ws.Range('E16').Value = '  +2.20%  '
$ws.Range('D17').Value = '26.470.75'
$ws.Range('E17').Value = '  +0.63%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '4.944'
$ws.Range('D18').Style = $cleanStyle
$ws.Range('E18').Value = '  +0.39%  '
$ws.Range('E19').Value = '  +0.15%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.99'
$ws.Range('D20').Style = $cleanStyle
$ws.Range('E20').Value = '  +1.23%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '191.13'
$ws.Range('D21').Style = $cleanStyle
$ws.Range('E21').Value = '  +0.91%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.256'
$ws.Range('D22').Style = $cleanStyle
$ws.Range('E22').Value = '  +0.92%  '
$ws.Range('E23').Value = '  +0.20%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '148.77'
$ws.Range('D24').Style = $cleanStyle
$ws.Range('E24').Value = '  +2.63%  '
$ws.Range('E25').Value = '  +6.05%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.914'
$ws.Range('D26').Style = $cleanStyle
$ws.Range('E26').Value = '  +2.44%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.79'
$ws.Range('D27').Style = $cleanStyle
$ws.Range('E27').Value = '  -0.07%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06220'
$ws.Range('D28').Style = $cleanStyle
$ws.Range('E28').Value = '  -5.75%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.388'
$ws.Range('D29').Style = $cleanStyle
$ws.Range('E29').Value = '  +2.29%  '
$ws.Range('E30').Value = '  +0.30%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.609'
$ws.Range('D31').Style = $cleanStyle
$ws.Range('E31').Value = '  +1.10%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.594'
$ws.Range('D32').Style = $cleanStyle
$ws.Range('E32').Value = '  +0.30%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.687'
$ws.Range('D33').Style = $cleanStyle
$ws.Range('E33').Value = '  +1.59%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.039'
$ws.Range('D34').Style = $cleanStyle
$ws.Range('E34').Value = '  +1.28%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6153'
$ws.Range('D35').Style = $cleanStyle
$ws.Range('E35').Value = '  -0.76%  '
$ws.Range('E36').Value = '  +0.51%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.757'
$ws.Range('D37').Style = $cleanStyle
$ws.Range('E37').Value = '  +1.39%  '
$ws.Range('E38').Value = '  +1.94%  '
$ws.Range('D39').Value = '1.118.57'
$ws.Range('E39').Value = '  +1.19%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.106'
$ws.Range('D40').Style = $cleanStyle
$ws.Range('E40').Value = '  -4.41%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8815'
$ws.Range('D41').Style = $cleanStyle
$ws.Range('E41').Value = '  +0.60%  '
$ws.Range('E42').Value = '  +0.16%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '101.20'
$ws.Range('D43').Style = $cleanStyle
$ws.Range('E43').Value = '  +0.74%  '
$ws.Range('D44').Value = '1.848.35'
$ws.Range('E44').Value = '  +1.15%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '57.61'
$ws.Range('D45').Style = $cleanStyle
$ws.Range('E45').Value = '  +1.39%  '
$ws.Range('E46').Value = '  -2.88%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.227'
$ws.Range('D47').Style = $cleanStyle
$ws.Range('E47').Value = '  +1.15%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.007'
$ws.Range('D48').Style = $cleanStyle
$ws.Range('E48').Value = '  -0.13%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05285'
$ws.Range('D49').Style = $cleanStyle
$ws.Range('E49').Value = '  +0.28%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.125'
$ws.Range('D50').Style = $cleanStyle
$ws.Range('E50').Value = '  +1.15%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4302'
$ws.Range('D51').Style = $cleanStyle
$ws.Range('E51').Value = '  +0.01%  '
